$wb = $excel.ActiveWorkbook
$ws46 = $wb.Worksheets.Item("EJ46")

# New column headers (shared strings J4601..J4606) replacing the generic placeholders
$ws46.Range("B1").Value = "J4601"
$ws46.Range("C1").Value = "J4602"
$ws46.Range("D1").Value = "J4603"
$ws46.Range("E1").Value = "J4604"
$ws46.Range("F1").Value = "J4605"
$ws46.Range("G1").Value = "J4606"

# Fill in the measured data (dx 0 .. 3)
$ws46.Range("B2").Value = 297.7
$ws46.Range("C2").Value = 219.8
$ws46.Range("D2").Value = 177.4
$ws46.Range("E2").Value = 132.8
$ws46.Range("F2").Value = 98.7
$ws46.Range("G2").Value = 74.2

$ws46.Range("B3").Value = 304.4
$ws46.Range("C3").Value = 225
$ws46.Range("D3").Value = 181.7
$ws46.Range("E3").Value = 135.7
$ws46.Range("F3").Value = 100.6
$ws46.Range("G3").Value = 75.7

$ws46.Range("B4").Value = 312.2
$ws46.Range("C4").Value = 230
$ws46.Range("D4").Value = 187.4
$ws46.Range("E4").Value = 138
$ws46.Range("F4").Value = 102.9
$ws46.Range("G4").Value = 77.6

$ws46.Range("B5").Value = 317.9
$ws46.Range("C5").Value = 235
$ws46.Range("D5").Value = 192.8
$ws46.Range("E5").Value = 141.1
$ws46.Range("F5").Value = 105.1
$ws46.Range("G5").Value = 79.4

$ws46.Range("B6").Value = 323.8
$ws46.Range("C6").Value = 240.6
$ws46.Range("D6").Value = 196.3
$ws46.Range("E6").Value = 143.5
$ws46.Range("F6").Value = 107.2
$ws46.Range("G6").Value = 81

$ws46.Range("B7").Value = 329.2
$ws46.Range("C7").Value = 245.9
$ws46.Range("D7").Value = 200.1
$ws46.Range("E7").Value = 146.5
$ws46.Range("F7").Value = 109.4
$ws46.Range("G7").Value = 82.8

$ws46.Range("B8").Value = 336.1
$ws46.Range("C8").Value = 251.6
$ws46.Range("D8").Value = 204.3
$ws46.Range("E8").Value = 149.1
$ws46.Range("F8").Value = 111.5
$ws46.Range("G8").Value = 84.3

# Move the active tab / selection from EJ43 to EJ46
$ws46.Activate()
$ws46.Range("G9").Select()
